$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2").Value = "2016-03-18 04:37:08"
$zhcn.Range("H2").Value = "2016-03-18 04:37:26"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2").Value = "2016-03-18 04:37:11"
$dede.Range("H2").Value = "2016-03-18 04:37:31"
